# TEST_CALENDAR_EXCEL_WEEK.xlsx — "dinamicamente aggiornare la data":
# shift every week-start date forward by one day and fix the day-of-month
# overflow bug (the old sheet kept counting past the end of March instead
# of rolling into April/May, e.g. "32/3/2019", "33/3/2019", ...).
# Also collapses the per-row B:H merges (which used to span 1370 rows) down
# to the four header merges actually used on row 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Row 2 holds the "W n" week-start dates (as plain D/M/Y text, one per
#    column from A2 to BB2). Recompute them as start_date + 1 day.
# ---------------------------------------------------------------------
$weekDates = [ordered]@{
    "A2"  = "23/3/2019";  "B2"  = "24/3/2019";  "C2"  = "25/3/2019";  "D2"  = "26/3/2019"
    "E2"  = "27/3/2019";  "F2"  = "28/3/2019";  "G2"  = "29/3/2019";  "H2"  = "30/3/2019"
    "I2"  = "31/3/2019";  "J2"  = "1/4/2019";   "K2"  = "2/4/2019";   "L2"  = "3/4/2019"
    "M2"  = "4/4/2019";   "N2"  = "5/4/2019";   "O2"  = "6/4/2019";   "P2"  = "7/4/2019"
    "Q2"  = "8/4/2019";   "R2"  = "9/4/2019";   "S2"  = "10/4/2019";  "T2"  = "11/4/2019"
    "U2"  = "12/4/2019";  "V2"  = "13/4/2019";  "W2"  = "14/4/2019";  "X2"  = "15/4/2019"
    "Y2"  = "16/4/2019";  "Z2"  = "17/4/2019";  "AA2" = "18/4/2019";  "AB2" = "19/4/2019"
    "AC2" = "20/4/2019";  "AD2" = "21/4/2019";  "AE2" = "22/4/2019";  "AF2" = "23/4/2019"
    "AG2" = "24/4/2019";  "AH2" = "25/4/2019";  "AI2" = "26/4/2019";  "AJ2" = "27/4/2019"
    "AK2" = "28/4/2019";  "AL2" = "29/4/2019";  "AM2" = "30/4/2019";  "AN2" = "1/5/2019"
    "AO2" = "2/5/2019";   "AP2" = "3/5/2019";   "AQ2" = "4/5/2019";   "AR2" = "5/5/2019"
    "AS2" = "6/5/2019";   "AT2" = "7/5/2019";   "AU2" = "8/5/2019";   "AV2" = "9/5/2019"
    "AW2" = "10/5/2019";  "AX2" = "11/5/2019";  "AY2" = "12/5/2019";  "AZ2" = "13/5/2019"
    "BA2" = "14/5/2019";  "BB2" = "15/5/2019"
}

# Stash the existing (rotated header) row-2 formatting in an unused cell far
# off the used range, so it can be restored after the writes below.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("ZZ1000").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Force the whole row to Text format first — otherwise entries such as
# "1/4/2019" get silently reinterpreted as an M/D/Y serial date instead of
# staying the literal D/M/Y text the sheet displays.
$ws.Range("A2:BB2").NumberFormat = "@"
foreach ($addr in $weekDates.Keys) {
    $ws.Range($addr).Value = $weekDates[$addr]
}

# Re-apply the original formatting (font / centered+rotated alignment) that
# the temporary Text number format above overwrote.
$ws.Range("ZZ1000").Copy() | Out-Null
$ws.Range("A2:BB2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("ZZ1000").Clear() | Out-Null

# ---------------------------------------------------------------------
# 2) The sheet used to merge B:H on every single data row (rows 10-19,
#    110-199, 1100-1370 — 371 merges total) even though only row 1 is
#    actually used. Drop all of those and merge the real header instead.
# ---------------------------------------------------------------------
$ws.Range("B10:H1370").UnMerge()

$ws.Range("B1:H1").Merge()
$ws.Range("I1:O1").Merge()
$ws.Range("P1:X1").Merge()
$ws.Range("Y1:AF1").Merge()
